$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the forecast table with a new row (row 39), carrying forward the
# same per-column formatting as the prior row (A38:E38) before writing the
# new values — this keeps column A's date style (s=2) without minting any
# unused/duplicate style entries.
$ws.Range("A38:E38").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.3398512689293476
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = 0.8571438361188566
